$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column A width adjustment (closest achievable value to 22.3125)
$ws.Columns.Item(1).ColumnWidth = 21.5

$ws.Cells.Item(2, 3).Value = -0.42
$ws.Cells.Item(2, 4).Value = -0.19
$ws.Cells.Item(2, 5).Value = -0.18
$ws.Cells.Item(2, 6).Value = -0.02
$ws.Cells.Item(2, 7).Value = -0.25
$ws.Cells.Item(2, 8).Value = -0.1
$ws.Cells.Item(2, 9).Value = -0.22
$ws.Cells.Item(2, 10).Value = -0.14
$ws.Cells.Item(2, 11).Value = -0.33
$ws.Cells.Item(2, 12).Value = -0.2
$ws.Cells.Item(2, 13).Value = -0.17

$ws.Cells.Item(3, 3).Value = -0.44
$ws.Cells.Item(3, 4).Value = -0.17
$ws.Cells.Item(3, 5).Value = -0.14
$ws.Cells.Item(3, 6).Value = 0.0
$ws.Cells.Item(3, 7).Value = -0.19
$ws.Cells.Item(3, 8).Value = -0.07
$ws.Cells.Item(3, 9).Value = -0.17
$ws.Cells.Item(3, 10).Value = -0.2
$ws.Cells.Item(3, 11).Value = -0.36
$ws.Cells.Item(3, 12).Value = -0.18
$ws.Cells.Item(3, 13).Value = -0.11

$ws.Cells.Item(4, 3).Value = -0.44
$ws.Cells.Item(4, 4).Value = -0.16
$ws.Cells.Item(4, 5).Value = -0.13
$ws.Cells.Item(4, 6).Value = 0.0
$ws.Cells.Item(4, 7).Value = -0.19
$ws.Cells.Item(4, 8).Value = -0.07
$ws.Cells.Item(4, 9).Value = -0.17
$ws.Cells.Item(4, 10).Value = -0.2
$ws.Cells.Item(4, 11).Value = -0.37
$ws.Cells.Item(4, 12).Value = -0.18
$ws.Cells.Item(4, 13).Value = -0.09

$ws.Cells.Item(5, 3).Value = -0.44
$ws.Cells.Item(5, 4).Value = -0.15
$ws.Cells.Item(5, 5).Value = -0.12
$ws.Cells.Item(5, 6).Value = -0.02
$ws.Cells.Item(5, 7).Value = -0.19
$ws.Cells.Item(5, 8).Value = -0.07
$ws.Cells.Item(5, 9).Value = -0.16
$ws.Cells.Item(5, 10).Value = -0.21
$ws.Cells.Item(5, 11).Value = -0.38
$ws.Cells.Item(5, 12).Value = -0.16
$ws.Cells.Item(5, 13).Value = -0.08

$ws.Cells.Item(6, 3).Value = -0.55
$ws.Cells.Item(6, 4).Value = -0.46
$ws.Cells.Item(6, 5).Value = -0.62
$ws.Cells.Item(6, 6).Value = -0.82
$ws.Cells.Item(6, 7).Value = -0.46
$ws.Cells.Item(6, 8).Value = 0.06
$ws.Cells.Item(6, 9).Value = -0.52
$ws.Cells.Item(6, 10).Value = -0.55
$ws.Cells.Item(6, 11).Value = -0.23
$ws.Cells.Item(6, 12).Value = -0.83
$ws.Cells.Item(6, 13).Value = -0.46

$ws.Cells.Item(7, 3).Value = -0.55
$ws.Cells.Item(7, 4).Value = -0.47
$ws.Cells.Item(7, 5).Value = -0.62
$ws.Cells.Item(7, 6).Value = -1.63
$ws.Cells.Item(7, 7).Value = -290.45
$ws.Cells.Item(7, 8).Value = 0.06
$ws.Cells.Item(7, 9).Value = -0.65
$ws.Cells.Item(7, 10).Value = -0.55
$ws.Cells.Item(7, 11).Value = -0.23
$ws.Cells.Item(7, 12).Value = -0.83
$ws.Cells.Item(7, 13).Value = -0.46

$ws.Cells.Item(8, 3).Value = -0.55
$ws.Cells.Item(8, 4).Value = -0.47
$ws.Cells.Item(8, 5).Value = -0.62
$ws.Cells.Item(8, 6).Value = -0.82
$ws.Cells.Item(8, 7).Value = -268.12
$ws.Cells.Item(8, 8).Value = 0.06
$ws.Cells.Item(8, 9).Value = -0.65
$ws.Cells.Item(8, 10).Value = -0.55
$ws.Cells.Item(8, 11).Value = -0.23
$ws.Cells.Item(8, 12).Value = -0.83
$ws.Cells.Item(8, 13).Value = -0.46

$ws.Cells.Item(9, 3).Value = -0.55
$ws.Cells.Item(9, 4).Value = -0.46
$ws.Cells.Item(9, 5).Value = -0.62
$ws.Cells.Item(9, 6).Value = -1.46
$ws.Cells.Item(9, 7).Value = -0.53
$ws.Cells.Item(9, 8).Value = 0.06
$ws.Cells.Item(9, 9).Value = -0.65
$ws.Cells.Item(9, 10).Value = -0.55
$ws.Cells.Item(9, 11).Value = -0.23
$ws.Cells.Item(9, 12).Value = -0.83
$ws.Cells.Item(9, 13).Value = -0.46

$ws.Cells.Item(10, 3).Value = -0.24
$ws.Cells.Item(10, 4).Value = 0.31
$ws.Cells.Item(10, 5).Value = -1.58
$ws.Cells.Item(10, 6).Value = 0.45
$ws.Cells.Item(10, 7).Value = 0.06
$ws.Cells.Item(10, 8).Value = 0.13
$ws.Cells.Item(10, 9).Value = 0.29
$ws.Cells.Item(10, 10).Value = 0.34
$ws.Cells.Item(10, 11).Value = -0.09
$ws.Cells.Item(10, 12).Value = -0.39
$ws.Cells.Item(10, 13).Value = 0.36

$ws.Cells.Item(11, 3).Value = -0.24
$ws.Cells.Item(11, 4).Value = 0.27
$ws.Cells.Item(11, 5).Value = -1.59
$ws.Cells.Item(11, 6).Value = 0.46
$ws.Cells.Item(11, 7).Value = 0.05
$ws.Cells.Item(11, 8).Value = 0.13
$ws.Cells.Item(11, 9).Value = 0.33
$ws.Cells.Item(11, 10).Value = 0.37
$ws.Cells.Item(11, 11).Value = -0.12
$ws.Cells.Item(11, 12).Value = -0.39
$ws.Cells.Item(11, 13).Value = 0.39

$ws.Cells.Item(12, 3).Value = -0.24
$ws.Cells.Item(12, 4).Value = 0.25
$ws.Cells.Item(12, 5).Value = -1.59
$ws.Cells.Item(12, 6).Value = 0.47
$ws.Cells.Item(12, 7).Value = 0.03
$ws.Cells.Item(12, 8).Value = 0.13
$ws.Cells.Item(12, 9).Value = 0.34
$ws.Cells.Item(12, 10).Value = 0.38
$ws.Cells.Item(12, 11).Value = -0.12
$ws.Cells.Item(12, 12).Value = -0.39
$ws.Cells.Item(12, 13).Value = 0.39

$ws.Cells.Item(13, 3).Value = -0.24
$ws.Cells.Item(13, 4).Value = 0.22
$ws.Cells.Item(13, 5).Value = -1.59
$ws.Cells.Item(13, 6).Value = 0.47
$ws.Cells.Item(13, 7).Value = -0.01
$ws.Cells.Item(13, 8).Value = 0.13
$ws.Cells.Item(13, 9).Value = 0.35
$ws.Cells.Item(13, 10).Value = 0.38
$ws.Cells.Item(13, 11).Value = -0.13
$ws.Cells.Item(13, 12).Value = -0.39
$ws.Cells.Item(13, 13).Value = 0.39

$ws.Cells.Item(14, 3).Value = -0.32
$ws.Cells.Item(14, 4).Value = 0.2
$ws.Cells.Item(14, 5).Value = -2.14
$ws.Cells.Item(14, 6).Value = -2.01
$ws.Cells.Item(14, 7).Value = -0.28
$ws.Cells.Item(14, 8).Value = 0.19
$ws.Cells.Item(14, 9).Value = -0.12
$ws.Cells.Item(14, 10).Value = -0.3
$ws.Cells.Item(14, 11).Value = -0.01
$ws.Cells.Item(14, 12).Value = -0.56
$ws.Cells.Item(14, 13).Value = -0.2

$ws.Cells.Item(15, 3).Value = -0.27
$ws.Cells.Item(15, 4).Value = 0.19
$ws.Cells.Item(15, 5).Value = -2.12
$ws.Cells.Item(15, 6).Value = -2.29
$ws.Cells.Item(15, 7).Value = -11883354.07
$ws.Cells.Item(15, 8).Value = -350819.04
$ws.Cells.Item(15, 9).Value = -0.12
$ws.Cells.Item(15, 10).Value = -0.3
$ws.Cells.Item(15, 11).Value = -0.01
$ws.Cells.Item(15, 12).Value = -0.56
$ws.Cells.Item(15, 13).Value = -0.2

$ws.Cells.Item(16, 3).Value = -0.27
$ws.Cells.Item(16, 4).Value = 0.18
$ws.Cells.Item(16, 5).Value = -2.06
$ws.Cells.Item(16, 6).Value = -2.6
$ws.Cells.Item(16, 7).Value = -52534948.8
$ws.Cells.Item(16, 8).Value = -1182358.04
$ws.Cells.Item(16, 9).Value = -0.12
$ws.Cells.Item(16, 10).Value = -0.3
$ws.Cells.Item(16, 11).Value = -0.01
$ws.Cells.Item(16, 12).Value = -0.56
$ws.Cells.Item(16, 13).Value = -0.2

$ws.Cells.Item(17, 3).Value = -0.27
$ws.Cells.Item(17, 4).Value = 0.17
$ws.Cells.Item(17, 5).Value = -2.14
$ws.Cells.Item(17, 6).Value = -3.4
$ws.Cells.Item(17, 7).Value = -315880791.11
$ws.Cells.Item(17, 8).Value = -3282813.07
$ws.Cells.Item(17, 9).Value = -0.12
$ws.Cells.Item(17, 10).Value = -0.3
$ws.Cells.Item(17, 11).Value = -0.01
$ws.Cells.Item(17, 12).Value = -0.56
$ws.Cells.Item(17, 13).Value = -0.19

$ws.Cells.Item(18, 3).Value = -0.24
$ws.Cells.Item(18, 4).Value = 0.31
$ws.Cells.Item(18, 5).Value = -0.12
$ws.Cells.Item(18, 6).Value = 0.47
$ws.Cells.Item(18, 7).Value = 0.06
$ws.Cells.Item(18, 8).Value = 0.19
$ws.Cells.Item(18, 9).Value = 0.35
$ws.Cells.Item(18, 10).Value = 0.38
$ws.Cells.Item(18, 11).Value = -0.01
$ws.Cells.Item(18, 12).Value = -0.16
$ws.Cells.Item(18, 13).Value = 0.39
